# Update "想去人数" (want-to-go count) values in F column for sheets
# "展览" (Exhibition) and "全部类型" (All types), as generated by the
# gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 333
$wsExhibit.Range("F3").Value = 8533
$wsExhibit.Range("F4").Value = 6246
$wsExhibit.Range("F5").Value = 549
$wsExhibit.Range("F8").Value = 73
$wsExhibit.Range("F9").Value = 335
$wsExhibit.Range("F10").Value = 1235
$wsExhibit.Range("F11").Value = 90

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 333
$wsAll.Range("F3").Value = 8533
$wsAll.Range("F4").Value = 6246
$wsAll.Range("F5").Value = 549
$wsAll.Range("F8").Value = 73
$wsAll.Range("F9").Value = 335
$wsAll.Range("F14").Value = 1235
$wsAll.Range("F15").Value = 90
